$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 13: 2020-06-12 (serial 43994) ---
$ws.Range("A13").Value = 11
$ws.Range("B13").Value = 43994
$ws.Range("C13").Value = 139196
$ws.Range("D13").Value = 197590
$ws.Range("E13").Value = 56928
$ws.Range("F13").Value = 16448
$ws.Range("G13").Value = 32.78

# Reuse existing formats: A column style (s=1) and the datetime style (s=2)
$ws.Range("A12").Copy() | Out-Null
$ws.Range("A13").PasteSpecial(-4122) | Out-Null
$ws.Range("B12").Copy() | Out-Null
$ws.Range("B13").PasteSpecial(-4122) | Out-Null

# --- Row 14: 2020-06-13 (serial 43995) ---
$ws.Range("A14").Value = 12
$ws.Range("B14").Value = 43995
$ws.Range("C14").Value = 142690
$ws.Range("D14").Value = 202139
$ws.Range("E14").Value = 56926
$ws.Range("F14").Value = 16872
$ws.Range("G14").Value = 32.66

# Reuse the A-column style for A14
$ws.Range("A12").Copy() | Out-Null
$ws.Range("A14").PasteSpecial(-4122) | Out-Null

# B14 gets a new date-only number format (YYYY-MM-DD) distinct from the
# existing datetime format used by B2:B13
$ws.Range("B14").NumberFormat = "YYYY-MM-DD"
